$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Faktúra Storno"
$ws.Range("B10").Value = "Faktúra Storno"

$ws.Range("A14").Select()
